# The commit adds two new price-record rows for "Lechuga" (Escarola variety)
# dated 2022-08-10 into the daily log sheet. In the original workbook these
# rows land at position 1191-1192, pushing all of the following rows down by
# two (old row 1191 becomes 1193, old row 1192 becomes 1194, ..., old row
# 1273 becomes 1275), which is exactly what the xml diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1191, shifting everything from the old row
# 1191 onward down by two rows.
$ws.Rows("1191:1192").Insert()

# New row 1191: Escarola / Primera, 2022-08-10 (serial 44783)
$row1191 = @(10, "Vega Modelo de Temuco", "La Araucanía", 44783, 9, 100112033, "Lechuga", "Escarola", "Primera", 500, 12000, 13000, 12600, "`$/caja 15 unidades", "Provincia del Elquí", 840, 15, "Hortaliza")
for ($i = 0; $i -lt $row1191.Length; $i++) {
    $ws.Cells.Item(1191, $i + 1).Value = $row1191[$i]
}

# New row 1192: Escarola / Segunda, 2022-08-10 (serial 44783)
$row1192 = @(10, "Vega Modelo de Temuco", "La Araucanía", 44783, 9, 100112033, "Lechuga", "Escarola", "Segunda", 100, 10000, 10000, 10000, "`$/caja 15 unidades", "Provincia del Elquí", 667, 15, "Hortaliza")
for ($i = 0; $i -lt $row1192.Length; $i++) {
    $ws.Cells.Item(1192, $i + 1).Value = $row1192[$i]
}
